$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("JMBG {jmbg}", $true, $false, $false, $false, $false, $true, 1, $false, "JMBG {jmbgNum}", 2)
